$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new time entry in row 14 (A14), matching the style/format of A13
$ws.Range("A14").Value = 0.068182870370370366
$ws.Range("A14").NumberFormat = $ws.Range("A13").NumberFormat

# Extend the SUM formulas to include the new row
$ws.Range("C2").Formula = "=SUM(A2:A14)"
$ws.Range("B3").Formula = "=SUM(A9:A14)"

# Update the active selection as in the target workbook
$ws.Range("E4").Select()
